$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in the Price (D) and Volume (E) columns
# are stored as text, matching the workbook's original inline-string data,
# rather than being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.129.52"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.893.78"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "245.96"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  +5.76%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "41.18"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("E9").Value = "  +4.38%  "
$ws.Range("D10").Value = "52.76"
$ws.Range("E10").Value = "  +12.66%  "
$ws.Range("D11").Value = "0.0716"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").Value = "0.0993"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "12.25"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "1.907.91"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "0.694"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "4.79"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "35.115.83"
$ws.Range("D19").Value = "71.94"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "0.0₃0816"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "240.06"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +25.23%  "
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "170.42"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "8.42"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "4.125.96"
$ws.Range("E31").Value = "  +20.85%  "
$ws.Range("D32").Value = "4.12"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.942"
$ws.Range("E33").Value = "  +15.48%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.0560"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "0.0637"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("D43").Value = "16.04"
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("D44").Value = "89.42"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.331.56"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("B46").Value = "MultiversX"
$ws.Range("C46").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D46").Value = "49.30"
$ws.Range("E46").Value = "  +40.75%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "2.77"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "6.46"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "2.077.44"
$ws.Range("E51").Value = "  +1.34%  "
